$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"1.863243333333333"
$ws.Range("H2").Value = [double]"5.589729999999999"
$ws.Range("I2").Value = [double]"0.6067417803684044"
$ws.Range("J2").Value = [double]"0.6067417803684044"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.03257366666666667"
$ws.Range("N2").Value = [double]"0.097721"
$ws.Range("O2").Value = [double]"0.001227793554179957"
$ws.Range("P2").Value = [double]"0.001227793554179957"
$ws.Range("Q2").Value = [double]"0.06069266725888888"
$ws.Range("R2").Value = [double]"0.54623400533"
$ws.Range("S2").Value = [double]"0.0007449536469879984"
$ws.Range("T2").Value = [double]"0.0007449536469879984"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"1.863243333333333"
$ws.Range("H3").Value = [double]"5.589729999999999"
$ws.Range("I3").Value = [double]"0.6067417803684044"
$ws.Range("J3").Value = [double]"0.6067417803684044"
$ws.Range("O3").Value = [double]"0.7662385783512358"
$ws.Range("P3").Value = [double]"0.7662385783512359"
$ws.Range("Q3").Value = [double]"37.87694023842332"
$ws.Range("R3").Value = [double]"340.8924621458099"
$ws.Range("S3").Value = [double]"0.464908959215784"
$ws.Range("T3").Value = [double]"0.464908959215784"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"1.863243333333333"
$ws.Range("H4").Value = [double]"5.589729999999999"
$ws.Range("I4").Value = [double]"0.6067417803684044"
$ws.Range("J4").Value = [double]"0.6067417803684044"
$ws.Range("M4").Value = [double]"6.169174666666667"
$ws.Range("N4").Value = [double]"18.507524"
$ws.Range("O4").Value = [double]"0.2325336280945842"
$ws.Range("P4").Value = [double]"0.2325336280945842"
$ws.Range("Q4").Value = [double]"11.49467356983555"
$ws.Range("R4").Value = [double]"103.45206212852"
$ws.Range("S4").Value = [double]"0.1410878675056324"
$ws.Range("T4").Value = [double]"0.1410878675056324"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.9841983333333334"
$ws.Range("H5").Value = [double]"2.952595"
$ws.Range("I5").Value = [double]"0.3204918210730839"
$ws.Range("J5").Value = [double]"0.3204918210730839"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.03257366666666667"
$ws.Range("N5").Value = [double]"0.097721"
$ws.Range("O5").Value = [double]"0.001227793554179957"
$ws.Range("P5").Value = [double]"0.001227793554179957"
$ws.Range("Q5").Value = [double]"0.03205894844388889"
$ws.Range("R5").Value = [double]"0.288530535995"
$ws.Range("S5").Value = [double]"0.0003934977920809287"
$ws.Range("T5").Value = [double]"0.0003934977920809286"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"0.9841983333333334"
$ws.Range("H6").Value = [double]"2.952595"
$ws.Range("I6").Value = [double]"0.3204918210730839"
$ws.Range("J6").Value = [double]"0.3204918210730839"
$ws.Range("O6").Value = [double]"0.7662385783512358"
$ws.Range("P6").Value = [double]"0.7662385783512359"
$ws.Range("Q6").Value = [double]"20.00727483496833"
$ws.Range("R6").Value = [double]"180.065473514715"
$ws.Range("S6").Value = [double]"0.2455731973522385"
$ws.Range("T6").Value = [double]"0.2455731973522385"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"0.9841983333333334"
$ws.Range("H7").Value = [double]"2.952595"
$ws.Range("I7").Value = [double]"0.3204918210730839"
$ws.Range("J7").Value = [double]"0.3204918210730839"
$ws.Range("M7").Value = [double]"6.169174666666667"
$ws.Range("N7").Value = [double]"18.507524"
$ws.Range("O7").Value = [double]"0.2325336280945842"
$ws.Range("P7").Value = [double]"0.2325336280945842"
$ws.Range("Q7").Value = [double]"6.071691424975556"
$ws.Range("R7").Value = [double]"54.64522282478"
$ws.Range("S7").Value = [double]"0.07452512592876452"
$ws.Range("T7").Value = [double]"0.07452512592876451"
$ws.Range("G8").Value = [double]"0.2234583333333333"
$ws.Range("H8").Value = [double]"0.6703750000000001"
$ws.Range("I8").Value = [double]"0.07276639855851162"
$ws.Range("J8").Value = [double]"0.07276639855851162"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.03257366666666667"
$ws.Range("N8").Value = [double]"0.097721"
$ws.Range("O8").Value = [double]"0.001227793554179957"
$ws.Range("P8").Value = [double]"0.001227793554179957"
$ws.Range("Q8").Value = [double]"0.00727885726388889"
$ws.Range("R8").Value = [double]"0.06550971537500001"
$ws.Range("S8").Value = [double]"8.934211511103031E-05"
$ws.Range("T8").Value = [double]"8.934211511103031E-05"
$ws.Range("G9").Value = [double]"0.2234583333333333"
$ws.Range("H9").Value = [double]"0.6703750000000001"
$ws.Range("I9").Value = [double]"0.07276639855851162"
$ws.Range("J9").Value = [double]"0.07276639855851162"
$ws.Range("O9").Value = [double]"0.7662385783512358"
$ws.Range("P9").Value = [double]"0.7662385783512359"
$ws.Range("Q9").Value = [double]"4.542572505708333"
$ws.Range("R9").Value = [double]"40.883152551375"
$ws.Range("S9").Value = [double]"0.05575642178321336"
$ws.Range("T9").Value = [double]"0.05575642178321336"
$ws.Range("G10").Value = [double]"0.2234583333333333"
$ws.Range("H10").Value = [double]"0.6703750000000001"
$ws.Range("I10").Value = [double]"0.07276639855851162"
$ws.Range("J10").Value = [double]"0.07276639855851162"
$ws.Range("M10").Value = [double]"6.169174666666667"
$ws.Range("N10").Value = [double]"18.507524"
$ws.Range("O10").Value = [double]"0.2325336280945842"
$ws.Range("P10").Value = [double]"0.2325336280945842"
$ws.Range("Q10").Value = [double]"1.378553489055556"
$ws.Range("R10").Value = [double]"12.4069814015"
$ws.Range("S10").Value = [double]"0.01692063466018723"
$ws.Range("T10").Value = [double]"0.01692063466018723"

Write-Host "Updated 120 cells"
